$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeekApr8")
$ws.Copy($ws)
